$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New problem rows: Complement of Base 10 Integer (62) & Number Complement (63) ---
# Populate in the same column-major order the original author used, so that
# the shared-string table is built up in the same order as the source edit.
$ws.Range("A62").Value = "Complement of Base 10 Integer"
$ws.Range("A63").Value = "Number Complement"
$ws.Range("G62").Value = "1009 - Complement of Base 10 Integer"
$ws.Range("G63").Value = "476 - Number Complement"

$ws.Range("B62").Value = "Binary"
$ws.Range("C62").Value = "No"
$ws.Range("D62").Value = "No"
$ws.Range("E62").Value = "Easy"
$ws.Range("F62").Value = "Easy"

$ws.Range("B63").Value = "Binary"
$ws.Range("C63").Value = "No"
$ws.Range("D63").Value = "No"
$ws.Range("E63").Value = "Easy"
$ws.Range("F63").Value = "Easy"

$ws.Hyperlinks.Add($ws.Range("G62"), "1009%20-%20Complement%20of%20Base%2010%20Integer")
$ws.Range("G62").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("G63"), "476%20-%20Number%20Complement")
$ws.Range("G63").Style = "Hyperlink"

# --- Extend data validation ranges to cover the new rows (2..63) ---
$ws.Range("E2:F63").Validation.Delete()
$ws.Range("E2:F63").Validation.Add(3, 1, 1, '"Easy, Medium, Hard"')

$ws.Range("C2:C63").Validation.Delete()
$ws.Range("C2:C63").Validation.Add(3, 1, 1, '"Yes, No"')
$ws.Range("C2:C63").Validation.IgnoreBlank = $false

$ws.Range("D2:D63").Validation.Delete()
$ws.Range("D2:D63").Validation.Add(3, 1, 1, '"Yes, No"')

$ws.Range("B2:B63").Validation.Delete()
$ws.Range("B2:B63").Validation.Add(3, 1, 1, '"Array, Binary, Dynamic Programming, Graph, Interval, Linked List, Matrix, String, Tree, Heap, Class Design"')

# --- Update view: scroll back to top and select O52, matching the saved view state ---
$ws.Range("A1").Select()
$ws.Range("O52").Select()
